$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2..41 in column E ("dias_defasagem") are incremented by 1 day,
# reflecting a later snapshot date used when recomputing the lag.
for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 5)
    $cell.Value = $cell.Value2 + 1
}

# Two records (rows 16 and 31) were actually refreshed on 23/10/2024,
# so their "ultima_atualizacao" date and lag (reset to 1 day) differ
# from the simple +1 pattern applied to the rest of the sheet.
$ws.Cells.Item(16, 4).Value = "23/10/2024"
$ws.Cells.Item(16, 5).Value = 1

$ws.Cells.Item(31, 4).Value = "23/10/2024"
$ws.Cells.Item(31, 5).Value = 1
